# recite_hooker.xlsx update
# - Append a new "review date" (2018-08-22 / serial 43334) column to rows 6-11
# - Add a new memorization entry in row 12 (text + first review date)
# - Update the sheet view's selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = 43334

# --- Helper: append the next date in the "staircase" of review dates for a
#     row by copying formatting from the last populated date cell in that
#     row (so it keeps numFmtId 14 / style index 4 instead of minting a new
#     custom number format), then writing the new serial date value.
# (Named parameters aren't reliably bound in this host, so pass positionally.)
function Add-ReviewDate($LastCell, $NextCell) {
    $ws.Range($LastCell).Copy()
    $ws.Range($NextCell).PasteSpecial(-4122) | Out-Null
    $ws.Range($NextCell).Value = $newDate
}

Add-ReviewDate "J6" "K6"
Add-ReviewDate "J7" "K7"
Add-ReviewDate "I8" "J8"
Add-ReviewDate "H9" "I9"
Add-ReviewDate "G10" "H10"
Add-ReviewDate "F11" "G11"

# --- Row 12: new memorization entry ---
$ws.Range("B12").Value = "Weigh what doth move the common sort so much to favour this innovation, and it shall soon appear unto you, that the force of particular reasons which for your several opinions are alleged is a thing whereof the multitude never did nor could so consider as to be therewith wholly carried; but certain general inducements are used to make them more saleable your cause in gross;"

# Give C12 the same date formatting as the other review-date cells, then set
# the first review date for this new entry.
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = $newDate

$ws.Rows.Item(12).RowHeight = 90

# --- Sheet view: scroll down slightly and move the active selection ---
$ws.Range("B13").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
